$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$rng = $ws.Range("A2:B1")
$rng = $ws.Range("B1,A2")
$rng.Borders.LineStyle = 1
$rng.Borders.Weight = 2
$rng.Font.Bold = $true
$rng.HorizontalAlignment = -4108
$rng.VerticalAlignment = -4160
